# Updated cryptos list on Fri Oct  4 21:56:48 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns, and replaces the
# VeChain row (51) with Stellar.
#
# For numeric-looking Price values, a leading apostrophe is used to force
# Excel to store them as text (matching the workbook's original inlineStr
# text cells, e.g. "556.26" rather than being auto-converted to the number
# 556.26). The Style is then reset to "Normal" so no stray number-format /
# quote-prefix styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.338.76'
$ws.Range("E2").Value = '  +2.55%  '
$ws.Range("D3").Value = '2.424.58'
$ws.Range("E3").Value = '  +3.21%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''556.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.19%  '
$ws.Range("D6").Value = '''143.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.77%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +1.86%  '
$ws.Range("D9").Value = '2.424.29'
$ws.Range("E9").Value = '  +3.27%  '
$ws.Range("E10").Value = '  +4.51%  '
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("E12").Value = '  +1.55%  '
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("D14").Value = '''26.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.28%  '
$ws.Range("E15").Value = '  +9.17%  '
$ws.Range("D16").Value = '2.862.30'
$ws.Range("E16").Value = '  +3.23%  '
$ws.Range("D17").Value = '62.169.14'
$ws.Range("E17").Value = '  +2.45%  '
$ws.Range("D18").Value = '2.426.42'
$ws.Range("E18").Value = '  +3.35%  '
$ws.Range("E19").Value = '  +4.18%  '
$ws.Range("D20").Value = '''4.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.98%  '
$ws.Range("D21").Value = '''324.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = '  +2.77%  '
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("D24").Value = '''1.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.20%  '
$ws.Range("D25").Value = '''64.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.40%  '
$ws.Range("D26").Value = '''9.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.33%  '
$ws.Range("D27").Value = '''572.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.78%  '
$ws.Range("E28").Value = '  +2.97%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = '''8.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.46%  '
$ws.Range("D31").Value = '0.0₃0938'
$ws.Range("E31").Value = '  +8.85%  '
$ws.Range("E32").Value = '  +5.84%  '
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("E35").Value = '  +4.62%  '
$ws.Range("D36").Value = '''5.76'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.73%  '
$ws.Range("E38").Value = '  +5.26%  '
$ws.Range("E39").Value = '  +2.52%  '
$ws.Range("E40").Value = '  +2.61%  '
$ws.Range("E41").Value = '  +1.70%  '
$ws.Range("D42").Value = '''149.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.04%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '''41.71'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.71%  '
$ws.Range("D45").Value = '''2.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +13.97%  '
$ws.Range("D46").Value = '''151.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.49%  '
$ws.Range("E47").Value = '  +2.32%  '
$ws.Range("D48").Value = '''0.0540'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.36%  '
$ws.Range("E49").Value = '  +6.66%  '
$ws.Range("E50").Value = '  +3.72%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '''0.0916'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.64%  '
